# "Generate Report for Archive" — localization-status report refresh.
#
# The status text "Ready for handoff" moves to "In Translation" everywhere
# it appears (it's a single shared string reused by the Overview summary
# row and both locale detail sheets), and the now-shorter text's column
# narrows to match on every sheet that shows it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: per-locale status columns E (zh-cn) and F (de-de) ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "In Translation"
$ovw.Range("F2").Value = "In Translation"
$ovw.Columns.Item(5).ColumnWidth = 12.42
$ovw.Columns.Item(6).ColumnWidth = 12.42

# --- zh-cn detail sheet: Status column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.42

# --- de-de detail sheet: Status column C ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.42
